$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column date for every existing data row
# (rows 2-99) from 45181 to 45182.
for ($r = 2; $r -le 99; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}

# Row 99 gains an explicit row height (15pt, custom height) in the new
# file - setting RowHeight reproduces the ht="15" customHeight="1" markup.
$ws.Rows.Item(99).RowHeight = 15

# Append the new record as row 100.
$ws.Cells.Item(100, 1).Value = "A 42496-2023"

$ws.Cells.Item(100, 2).Value = 45180
$ws.Cells.Item(100, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(100, 3).Value = 45182
$ws.Cells.Item(100, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(100, 4).Value = "NORRBOTTENS LÄN"
$ws.Cells.Item(100, 5).Value = "KIRUNA"

$ws.Cells.Item(100, 7).Value = 9.300000000000001
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).Value = 0
$ws.Cells.Item(100, 14).Value = 0
$ws.Cells.Item(100, 15).Value = 0
$ws.Cells.Item(100, 16).Value = 0
$ws.Cells.Item(100, 17).Value = 0

# R100 mirrors R99's empty, wrap-text-styled cell.
$ws.Cells.Item(100, 18).WrapText = $true
